# Created configure template for frequency.
#
# Target sheet is "Translation" (the active sheet / second tab of the
# workbook). It holds a simple B:F table of translation text entries:
#   B = Text Id, C = Typography, D = Alignment, E = Direction, F = Label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo in the existing "Clock" row (row 8): "CLOSK" -> "CLOCK"
$ws.Range("F8").Value = "CLOCK"

# Duplicate the "Session setup" row as a template for five new rows,
# configuring each for a different menu entry (frequency measurement,
# input, clock, MES setup, session setup).
$ws.Range("B11").Value = "SingleUseId8"
$ws.Range("C11").Value = "Default"
$ws.Range("D11").Value = "Left"
$ws.Range("E11").Value = "LTR"
$ws.Range("F11").Value = "FREQUENCY MEASUREMENT"

$ws.Range("B12").Value = "SingleUseId9"
$ws.Range("C12").Value = "Default"
$ws.Range("D12").Value = "Left"
$ws.Range("E12").Value = "LTR"
$ws.Range("F12").Value = "INPUT"

$ws.Range("B13").Value = "SingleUseId10"
$ws.Range("C13").Value = "Default"
$ws.Range("D13").Value = "Left"
$ws.Range("E13").Value = "LTR"
$ws.Range("F13").Value = "CLOck"

$ws.Range("B14").Value = "SingleUseId11"
$ws.Range("C14").Value = "Default"
$ws.Range("D14").Value = "Left"
$ws.Range("E14").Value = "LTR"
$ws.Range("F14").Value = "MES SETUP"

$ws.Range("B15").Value = "SingleUseId12"
$ws.Range("C15").Value = "Default"
$ws.Range("D15").Value = "Left"
$ws.Range("E15").Value = "LTR"
$ws.Range("F15").Value = "SESSION SETUP"
